$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2337.4092
$ws.Range("I112").Value = 1083.2
$ws.Range("J112").Value = 2706.2942
$ws.Range("K112").Value = 3249.6
$ws.Range("L112").Value = 8118.882599999999
$ws.Range("M112").Value = -2141.6
$ws.Range("N112").Value = -10334.8826

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 11405
$ws.Range("I113").Value = 8503
$ws.Range("J113").Value = 20111
$ws.Range("K113").Value = 8503
$ws.Range("L113").Value = 20111
$ws.Range("M113").Value = -5249
$ws.Range("N113").Value = -26619

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2022.1111
$ws.Range("I127").Value = 1171.5714
$ws.Range("J127").Value = 4999
$ws.Range("K127").Value = 3514.7142
$ws.Range("L127").Value = 14997
$ws.Range("M127").Value = 1445.2858
$ws.Range("N127").Value = -24917

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5688.5557
$ws.Range("I132").Value = 5739.5
$ws.Range("J132").Value = 5484.778
$ws.Range("K132").Value = 17218.5
$ws.Range("L132").Value = 16454.334
$ws.Range("M132").Value = -14688.5
$ws.Range("N132").Value = -21514.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9604.733
$ws.Range("I137").Value = 12699.35
$ws.Range("J137").Value = 3415.5
$ws.Range("K137").Value = 38098.05
$ws.Range("L137").Value = 10246.5
$ws.Range("M137").Value = -35548.05
$ws.Range("N137").Value = -15346.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2329.6365
$ws.Range("I138").Value = 1422.9062
$ws.Range("J138").Value = 3591.1738
$ws.Range("K138").Value = 4268.7186
$ws.Range("L138").Value = 10773.5214
$ws.Range("M138").Value = 871.2813999999998
$ws.Range("N138").Value = -21053.5214

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3473.175
$ws.Range("I141").Value = 3322.5557
$ws.Range("J141").Value = 4828.75
$ws.Range("K141").Value = 9967.667099999999
$ws.Range("L141").Value = 14486.25
$ws.Range("M141").Value = -4787.667099999999
$ws.Range("N141").Value = -24846.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 21033.334
$ws.Range("I61").Value = 29000
$ws.Range("J61").Value = 5100
$ws.Range("K61").Value = 29000
$ws.Range("L61").Value = 5100
$ws.Range("M61").Value = -28788
$ws.Range("N61").Value = -5524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6960.4165
$ws.Range("I74").Value = 1627.75
$ws.Range("J74").Value = 9626.75
$ws.Range("K74").Value = 1627.75
$ws.Range("L74").Value = 9626.75
$ws.Range("M74").Value = -753.75
$ws.Range("N74").Value = -11374.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6960.4165
$ws.Range("I77").Value = 1627.75
$ws.Range("J77").Value = 9626.75
$ws.Range("K77").Value = 8138.75
$ws.Range("L77").Value = 48133.75
$ws.Range("M77").Value = -3770.75
$ws.Range("N77").Value = -56869.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1298712.2
$ws.Range("I122").Value = 3613.9546
$ws.Range("J122").Value = 3673059
$ws.Range("K122").Value = 10841.8638
$ws.Range("L122").Value = 11019177
$ws.Range("M122").Value = -8391.863799999999
$ws.Range("N122").Value = -11024077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 40433.332
$ws.Range("I131").Value = 40650
$ws.Range("J131").Value = 40000
$ws.Range("K131").Value = 40650
$ws.Range("L131").Value = 40000
$ws.Range("M131").Value = -35610
$ws.Range("N131").Value = -50080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3573.8374
$ws.Range("I132").Value = 3352.647
$ws.Range("J132").Value = 4827.25
$ws.Range("K132").Value = 10057.941
$ws.Range("L132").Value = 14481.75
$ws.Range("M132").Value = -7527.940999999999
$ws.Range("N132").Value = -19541.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 21033.334
$ws.Range("I136").Value = 29000
$ws.Range("J136").Value = 5100
$ws.Range("K136").Value = 87000
$ws.Range("L136").Value = 15300
$ws.Range("M136").Value = -84450
$ws.Range("N136").Value = -20400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 52000
$ws.Range("I57").Value = 52000
$ws.Range("K57").Value = 52000
$ws.Range("M57").Value = -51280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 14750
$ws.Range("J74").Value = 14750
$ws.Range("L74").Value = 14750
$ws.Range("N74").Value = -16622

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value = 14750
$ws.Range("J77").Value = 14750
$ws.Range("L77").Value = 44250
$ws.Range("N77").Value = -53610

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22628.75
$ws.Range("I134").Value = 22628.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 67886.25
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -65351.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 52000
$ws.Range("I136").Value = 52000
$ws.Range("K136").Value = 52000
$ws.Range("M136").Value = -46900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 22222
$ws.Range("J131").Value = 22222
$ws.Range("L131").Value = 22222
$ws.Range("N131").Value = -32302

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1472.3334
$ws.Range("I134").Value = 1020.6923
$ws.Range("K134").Value = 3062.0769
$ws.Range("M134").Value = -527.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 289.7143
$ws.Range("I92").Value = 259.5
$ws.Range("J92").Value = 330
$ws.Range("K92").Value = 778.5
$ws.Range("L92").Value = 990
$ws.Range("M92").Value = 469.5
$ws.Range("N92").Value = -3486

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 40247.332
$ws.Range("I97").Value = 75047.5
$ws.Range("J97").Value = 475.7143
$ws.Range("K97").Value = 225142.5
$ws.Range("L97").Value = 1427.1429
$ws.Range("M97").Value = -224646.5
$ws.Range("N97").Value = -2419.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 913.8570999999999
$ws.Range("I98").Value = 1022.8333
$ws.Range("J98").Value = 832.125
$ws.Range("K98").Value = 3068.4999
$ws.Range("L98").Value = 2496.375
$ws.Range("M98").Value = -1570.4999
$ws.Range("N98").Value = -5492.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 11401.667
$ws.Range("I130").Value = 3340.6667
$ws.Range("J130").Value = 15432.167
$ws.Range("K130").Value = 10022.0001
$ws.Range("L130").Value = 46296.501
$ws.Range("M130").Value = -5002.000100000001
$ws.Range("N130").Value = -56336.501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2272.5833
$ws.Range("I132").Value = 2388.3635
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 7165.0905
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -4635.0905
$ws.Range("N132").Value = -8057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 416889.56
$ws.Range("I132").Value = 553149.4399999999
$ws.Range("J132").Value = 8110
$ws.Range("K132").Value = 1659448.32
$ws.Range("L132").Value = 24330
$ws.Range("M132").Value = -1656918.32
$ws.Range("N132").Value = -29390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 7745465.5
$ws.Range("I5").Value = 17333332
$ws.Range("J5").Value = 554565
$ws.Range("K5").Value = 17333332
$ws.Range("L5").Value = 554565
$ws.Range("M5").Value = -17333220
$ws.Range("N5").Value = -554789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 27199.75
$ws.Range("I81").Value = 50400
$ws.Range("J81").Value = 3999.5
$ws.Range("K81").Value = 100800
$ws.Range("L81").Value = 7999
$ws.Range("M81").Value = -99739
$ws.Range("N81").Value = -10121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 27199.75
$ws.Range("I84").Value = 50400
$ws.Range("J84").Value = 3999.5
$ws.Range("K84").Value = 504000
$ws.Range("L84").Value = 39995
$ws.Range("M84").Value = -498696
$ws.Range("N84").Value = -50603
